{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph is the existing \"Introduction\" paragraph.\nconst introPara = body.paragraphs.items[0];\nintroPara.alignment = Word.Alignment.justified;\n\n// Insert a new paragraph after it containing the body text, also justified.\nconst newPara = introPara.insertParagraph(\n  \"A public school district consists of one or more public schools operated under the supervision of an elected or appointed school committee and a superintendent. The majority of school districts serve a single city or town, and are considered a department of the municipal government. Two or more municipalities can also join together to form a regional school district, which is considered a separate and independent unit of local government. A regional school district can offer all grades (preK-12), just certain grades (for example, just elementary grades or just high school), or just certain types of instruction (for example, vocational and technical programs).      Number of School Districts in this state are 399, Number of Schools 1,827, Enrollment 914,959, and Grades Served PK-12 (dese, 2024). The annual high school dropout statistics represents a snapshot of those students who dropped out of school in any given year.  This study will present the school district dropout rate in Massachusetts.\",\n  Word.InsertLocation.after\n);\nnewPara.alignment = Word.Alignment.justified;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Justify the existing \"Introduction\" paragraph.\n$p1 = $d.Paragraphs.First\n$p1.Alignment = 3  # wdAlignParagraphJustify\n\n# Insert a new paragraph after it with the body text, also justified.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = 'A public school district consists of one or more public schools operated under the supervision of an elected or appointed school committee and a superintendent. The majority of school districts serve a single city or town, and are considered a department of the municipal government. Two or more municipalities can also join together to form a regional school district, which is considered a separate and independent unit of local government. A regional school district can offer all grades (preK-12), just certain grades (for example, just elementary grades or just high school), or just certain types of instruction (for example, vocational and technical programs).      Number of School Districts in this state are 399, Number of Schools 1,827, Enrollment 914,959, and Grades Served PK-12 (dese, 2024). The annual high school dropout statistics represents a snapshot of those students who dropped out of school in any given year.  This study will present the school district dropout rate in Massachusetts.'\n$p2.Alignment = 3  # wdAlignParagraphJustify\n"}
